{"js": "// Map of old text -> new text, derived from the authoritative diff.\nconst replacements = [\n  [\"2024-09-20 Friday\", \"2024-09-21 Saturday\"],\n  [\"581\u00f77=83, 0\", \"340\u00f74=85, 0\"],\n  [\"263\u00f79=29, 2\", \"863\u00f76=143, 5\"],\n  [\"958\u00f77=136, 6\", \"226\u00f72=113, 0\"],\n  [\"848\u00f77=121, 1\", \"934\u00f74=233, 2\"],\n  [\"440\u00f72=220, 0\", \"105\u00f72=52, 1\"],\n  [\"248\u00f77=35, 3\", \"256\u00f76=42, 4\"],\n  [\"783\u00f73=261, 0\", \"375\u00f79=41, 6\"],\n  [\"637\u00f79=70, 7\", \"220\u00f78=27, 4\"],\n  [\"288\u00f79=32, 0\", \"450\u00f72=225, 0\"],\n  [\"291\u00f75=58, 1\", \"822\u00f76=137, 0\"],\n  [\"821\u00f75=164, 1\", \"544\u00f76=90, 4\"],\n  [\"484\u00f74=121, 0\", \"925\u00f72=462, 1\"],\n  [\"416\u00f73=138, 2\", \"979\u00f78=122, 3\"],\n  [\"433\u00f77=61, 6\", \"756\u00f72=378, 0\"],\n  [\"112\u00f73=37, 1\", \"853\u00f79=94, 7\"],\n  [\"525\u00f79=58, 3\", \"531\u00f75=106, 1\"],\n  [\"332\u00f77=47, 3\", \"524\u00f73=174, 2\"],\n  [\"843\u00f72=421, 1\", \"395\u00f77=56, 3\"],\n  [\"661\u00f77=94, 3\", \"829\u00f72=414, 1\"],\n  [\"722\u00f77=103, 1\", \"352\u00f76=58, 4\"],\n  [\"318\u00f74=79, 2\", \"463\u00f73=154, 1\"],\n  [\"525\u00f74=131, 1\", \"703\u00f78=87, 7\"],\n  [\"201\u00f76=33, 3\", \"411\u00f72=205, 1\"],\n  [\"608\u00f76=101, 2\", \"222\u00f72=111, 0\"],\n  [\"660\u00f72=330, 0\", \"167\u00f79=18, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Map of old text -> new text, derived from the authoritative diff.\n$pairs = @(\n  @(\"2024-09-20 Friday\", \"2024-09-21 Saturday\"),\n  @(\"581\u00f77=83, 0\", \"340\u00f74=85, 0\"),\n  @(\"263\u00f79=29, 2\", \"863\u00f76=143, 5\"),\n  @(\"958\u00f77=136, 6\", \"226\u00f72=113, 0\"),\n  @(\"848\u00f77=121, 1\", \"934\u00f74=233, 2\"),\n  @(\"440\u00f72=220, 0\", \"105\u00f72=52, 1\"),\n  @(\"248\u00f77=35, 3\", \"256\u00f76=42, 4\"),\n  @(\"783\u00f73=261, 0\", \"375\u00f79=41, 6\"),\n  @(\"637\u00f79=70, 7\", \"220\u00f78=27, 4\"),\n  @(\"288\u00f79=32, 0\", \"450\u00f72=225, 0\"),\n  @(\"291\u00f75=58, 1\", \"822\u00f76=137, 0\"),\n  @(\"821\u00f75=164, 1\", \"544\u00f76=90, 4\"),\n  @(\"484\u00f74=121, 0\", \"925\u00f72=462, 1\"),\n  @(\"416\u00f73=138, 2\", \"979\u00f78=122, 3\"),\n  @(\"433\u00f77=61, 6\", \"756\u00f72=378, 0\"),\n  @(\"112\u00f73=37, 1\", \"853\u00f79=94, 7\"),\n  @(\"525\u00f79=58, 3\", \"531\u00f75=106, 1\"),\n  @(\"332\u00f77=47, 3\", \"524\u00f73=174, 2\"),\n  @(\"843\u00f72=421, 1\", \"395\u00f77=56, 3\"),\n  @(\"661\u00f77=94, 3\", \"829\u00f72=414, 1\"),\n  @(\"722\u00f77=103, 1\", \"352\u00f76=58, 4\"),\n  @(\"318\u00f74=79, 2\", \"463\u00f73=154, 1\"),\n  @(\"525\u00f74=131, 1\", \"703\u00f78=87, 7\"),\n  @(\"201\u00f76=33, 3\", \"411\u00f72=205, 1\"),\n  @(\"608\u00f76=101, 2\", \"222\u00f72=111, 0\"),\n  @(\"660\u00f72=330, 0\", \"167\u00f79=18, 5\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
